$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Valor Mora" (F) and "Salario Basico" (G) values between the
# row for ALEXI WALTER PAJARO ACOSTA (row 18) and the row for
# LAURA VANESSA BLANCO VELASQUEZ (row 16), since the database values had
# been entered in the wrong rows.

$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2
$f18 = $ws.Range("F18").Value2
$g18 = $ws.Range("G18").Value2

$ws.Range("F16").Value2 = $f18
$ws.Range("G16").Value2 = $g18
$ws.Range("F18").Value2 = $f16
$ws.Range("G18").Value2 = $g16
